$d = $word.ActiveDocument

# Locate the target paragraph: style "a3" (Основной) whose whole text is
# exactly "В " followed by the paragraph mark (it also holds the lone
# "_GoBack" bookmark at its end). This is the paragraph that the commit
# extends with new sentences and then splits in two.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "В `r" -and $para.Style.NameLocal -eq "Основной") {
        $target = $para
    }
}

if ($target -ne $null) {
    $targetRange = $target.Range

    # Build the two replacement paragraphs (same "a3" style, no paragraph
    # mark language override) with the expanded wording. The bookmark that
    # used to sit right after "В " now ends up at the end of the second
    # paragraph, and the spell-checker marks around "мокапов"/"Мокапы" are
    # reproduced with explicit proofErr markers, just like Word would add
    # them.
    $innerXml = (
        '<w:p><w:pPr><w:pStyle w:val="a3"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">В </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">рамках разработки оконного приложения создан интерфейс пользователя в виде </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>мокапов</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>. Эти визуальные представления наглядно показывают структуру приложения, его основные элементы и функциональность.</w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="a3"/></w:pPr>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Мокапы</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> интерфейса, такие как </w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '</w:p>'
    )

    $xmlFrag = (
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $innerXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    )

    $null = $targetRange.InsertXML($xmlFrag)
}

Write-Output "done"
